$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.479.03'
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  +1.56%  '
$style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.000.96'
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  +4.43%  '
$ws.Range('E4').Value = '  -0.01%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '323.99'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('E6').Value = '  -0.04%  '
$style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5105'
$ws.Range('D7').Style = $style
$ws.Range('E7').Value = '  +1.46%  '
$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4193'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  +4.24%  '
$style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08734'
$ws.Range('D9').Style = $style
$ws.Range('E9').Value = '  +5.90%  '
$ws.Range('E10').Value = '  +2.17%  '
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '43.04'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  +2.34%  '
$ws.Range('E12').Value = '  +4.47%  '
$style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.996.34'
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  +4.23%  '
$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.580'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  +2.44%  '
$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.440'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  +1.88%  '
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  -0.10%  '
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '94.32'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  +2.32%  '
$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001115'
$ws.Range('D18').Style = $style
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06477'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  -0.20%  '
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.92'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  +4.01%  '
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.187'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  +4.02%  '
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '30.543.10'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  +1.63%  '
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.83'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  +4.94%  '
$ws.Range('E25').Value = '  +1.48%  '
$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.229.91'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  +4.36%  '
$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.30'
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  +0.25%  '
$style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '163.38'
$ws.Range('D28').Style = $style
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('E29').Value = '  +4.37%  '
$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '131.75'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  +2.18%  '
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('E32').Value = '  +1.05%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.069'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  +0.91%  '
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.853'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  +0.73%  '
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.334'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  +10.35%  '
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02522'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  +3.22%  '
$style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.444'
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  +1.46%  '
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06603'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  +2.77%  '
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '12.46'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  +9.54%  '
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2200'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  +1.67%  '
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.047'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  +1.07%  '
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.6627'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  +3.20%  '
$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.235'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  +1.52%  '
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.65'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  +2.53%  '
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6167'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  +2.77%  '
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.202'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +1.62%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.666'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  +0.81%  '
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.273'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  +4.74%  '
$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '124.44'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  +1.07%  '
$style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '80.34'
$ws.Range('D50').Style = $style
$ws.Range('E50').Value = '  +1.87%  '
$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06902'
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  +1.57%  '
